# Add login verification and globalExceptionHandler
#
# Functional change recorded in the XML diff: a new worksheet "report1" is
# added (after "Sheet1"), populated with a small report table, a thin box
# border style is introduced and applied to part of that table, and the
# selection/active-sheet state moves to the new sheet. Sheet1's own
# selection is left parked at D15.

$wb = $excel.ActiveWorkbook

# --- Sheet1: update the lingering selection before we move away from it ---
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("D15").Select()

# --- Add the new sheet right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "report1"

# Header row
$headers = @("a", "b", "c", "d", "e", "f")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2: two numbers followed by four text values
$ws2.Cells.Item(2, 1).Value = 1
$ws2.Cells.Item(2, 2).Value = 2
$ws2.Cells.Item(2, 3).Value = "你好呀"
$ws2.Cells.Item(2, 4).Value = "hello world"
$ws2.Cells.Item(2, 5).Value = "xxx"
$ws2.Cells.Item(2, 6).Value = "abc"

# Rows 3-6: repeating 1,2,3,4,5,6
for ($r = 3; $r -le 6; $r++) {
    $ws2.Cells.Item($r, 1).Value = 1
    $ws2.Cells.Item($r, 2).Value = 2
    $ws2.Cells.Item($r, 3).Value = 3
    $ws2.Cells.Item($r, 4).Value = 4
    $ws2.Cells.Item($r, 5).Value = 5
    $ws2.Cells.Item($r, 6).Value = 6
}

# Thin box border around C2:F6
$ws2.Range("C2:F6").Borders.LineStyle = 1

# Leave the new sheet active with C3 selected
[void]$ws2.Range("C3").Select()
